$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.732.04"
$ws.Range("E2").Value = "  -0.75%  "

$ws.Range("D3").Value = "2.033.03"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'228.02"
$ws.Range("E5").Value = "  -0.31%  "

$ws.Range("E6").Value = "  -1.54%  "

$ws.Range("D7").Value = "'60.18"
$ws.Range("E7").Value = "  -0.57%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.377"
$ws.Range("E9").Value = "  -2.00%  "

$ws.Range("E10").Value = "  +2.64%  "

$ws.Range("E11").Value = "  +0.38%  "

$ws.Range("D12").Value = "'14.65"
$ws.Range("E12").Value = "  -0.46%  "

$ws.Range("D13").Value = "2.334.52"
$ws.Range("E13").Value = "  -0.78%  "

$ws.Range("D14").Value = "'21.05"
$ws.Range("E14").Value = "  +0.10%  "

$ws.Range("E15").Value = "  +2.16%  "

$ws.Range("D16").Value = "'5.18"
$ws.Range("E16").Value = "  -2.65%  "

$ws.Range("D17").Value = "2.022.16"
$ws.Range("E17").Value = "  -1.79%  "

$ws.Range("D18").Value = "37.697.63"
$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("D19").Value = "'69.59"
$ws.Range("E19").Value = "  -0.11%  "

$ws.Range("D20").Value = "'5.88"
$ws.Range("E20").Value = "  -6.09%  "

$ws.Range("E21").Value = "  -0.72%  "

$ws.Range("D22").Value = "'223.78"
$ws.Range("E22").Value = "  -0.79%  "

$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("E24").Value = "  -0.81%  "

$ws.Range("E25").Value = "  +3.03%  "

$ws.Range("E26").Value = "  +1.89%  "

$ws.Range("D27").Value = "'167.41"
$ws.Range("E27").Value = "  +1.16%  "

$ws.Range("E28").Value = "  -2.49%  "

$ws.Range("D29").Value = "'18.77"
$ws.Range("E29").Value = "  -1.18%  "

$ws.Range("D30").Value = "'1.26"
$ws.Range("E30").Value = "  -2.52%  "

$ws.Range("E31").Value = "  +0.47%  "

$ws.Range("E32").Value = "  +9.12%  "

$ws.Range("E33").Value = "  -2.93%  "

$ws.Range("E34").Value = "  +0.37%  "

$ws.Range("E35").Value = "  -1.06%  "

$ws.Range("D36").Value = "'6.43"
$ws.Range("E36").Value = "  +3.03%  "

$ws.Range("E37").Value = "  +1.04%  "

$ws.Range("D38").Value = "'3.42"
$ws.Range("E38").Value = "  +4.56%  "

$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("D40").Value = "'18.14"
$ws.Range("E40").Value = "  +7.57%  "

$ws.Range("D41").Value = "1.536.93"
$ws.Range("E41").Value = "  +1.09%  "

$ws.Range("E42").Value = "  +0.32%  "

$ws.Range("D43").Value = "'96.05"
$ws.Range("E43").Value = "  -1.24%  "

$ws.Range("E44").Value = "  -2.40%  "

$ws.Range("D45").Value = "'0.0911"
$ws.Range("E45").Value = "  -1.20%  "

$ws.Range("E46").Value = "  -1.61%  "

$ws.Range("E47").Value = "  -0.53%  "

$ws.Range("E48").Value = "  -0.39%  "

$ws.Range("E49").Value = "  +0.24%  "

$ws.Range("D50").Value = "'7.10"
$ws.Range("E50").Value = "  +1.37%  "

$ws.Range("D51").Value = "2.223.82"
$ws.Range("E51").Value = "  -0.84%  "
